# Completes the PBCoreDigitalInstantiation upload Mapper fixture:
#   - shifts the aapb_preservation_lto / aapb_preservation_disk columns one
#     column to the right (F->G, and a new F)
#   - adds a new "DigitalInstantiation.location" column in the old E slot
#     (with sample value "Master")
#   - the old "DigitalInstantiation.holding_institution" header in D1 is
#     overwritten with "DigitalInstantiation.generations" (mirroring C1)
#     while its former value moves down into D2 as "Master" after shift
#
# Writes happen right-to-left / value-before-header so the shared-string
# table ends up appended in the same order Excel produced it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the two trailing columns right by one (G gets old F's value, F gets
# old E's value).
$ws.Range("G1").Value = "DigitalInstantition.aapb_preservation_disk"
$ws.Range("F1").Value = "DigitalInstantition.aapb_preservation_lto"
$ws.Range("D1").Value = "DigitalInstantiation.generations"

$ws.Range("G2").Value = "disky mc diskface"
$ws.Range("F2").Value = "fhqwhgads"
$ws.Range("E2").Value = "American Archive of Public Broadcasting"
$ws.Range("D2").Value = "Master"

# New column: DigitalInstantiation.location
$ws.Range("E1").Value = "DigitalInstantiation.location"

$ws.Range("E3").Select() | Out-Null
